# Bugfixed the naive forecaster component module
# Update the YoY forecast vectors in columns C (values for prior-year pair)
# and E (values for forecast-year pair) with corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "E2"  = 1.516248937663556
    "C3"  = 1.514319819128396
    "E3"  = 1.602279001294704
    "C4"  = 2.007652128026982
    "E4"  = 1.768040115052738
    "C5"  = 1.022680528298392
    "E5"  = 1.510468690286459
    "C6"  = 0.9070039918702477
    "E6"  = 1.042579621507111
    "C7"  = 1.022042907336096
    "E7"  = 1.083482333436536
    "C8"  = 1.34489417553354
    "E8"  = 1.095916825800991
    "C9"  = 1.277042522796856
    "E9"  = 1.293136192195643
    "C10" = 2.615369162917314
    "E10" = 1.825134644920934
    "C11" = 2.618053282882693
    "E11" = 2.406099663413808
    "C12" = 1.431088640641853
    "E12" = 2.21629047761287
    "C13" = 2.136062314641141
    "E13" = 2.031764787322499
    "C14" = 2.166968775134936
    "E14" = 1.984987808509886
    "C15" = 2.459440348120401
    "E15" = 2.526389380645511
    "C16" = 0.8171929556848756
    "E16" = 1.509741350988136
    "C17" = 0.8766015904249524
    "E17" = 2.477445663648559
    "C18" = 1.554086551645839
    "E18" = 0.9888012784191602
    "C19" = 0.5837948599211717
    "E19" = 1.328924132093245
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
